$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("F2").Value = 0.258
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = -1261.9
$ws.Range("L2").Value = -0.2580467056562104
$ws.Range("M2").Value = 21.17
$ws.Range("N2").Value = 0.00272107969151671
$ws.Range("O2").Value = -0.01677628972184801
$ws.Range("R2").Value = -0
$ws.Range("S2").Value = 21.17
$ws.Range("U2").Value = 30827.3
$ws.Range("V2").Value = 3.962377892030848
$ws.Range("W2").Value = -0.03168289030325721
$ws.Range("X2").Value = 0.2490925801383365
$ws.Range("Y2").Value = -0.2807754704415937
$ws.Range("Z2").Value = 0.0787071313148221
$ws.Range("AA2").Value = 0
$ws.Range("AB2").Value = 0.08516646298621965
$ws.Range("AC2").Value = -0.08516646298621965
$ws.Range("AD2").Value = 54704.2
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 54704.2
$ws.Range("AG2").Value = 23876.9
$ws.Range("AH2").Value = 0.8754885234987405
$ws.Range("AI2").Value = 0.6300476242578504
$ws.Range("AJ2").Value = 0.7542399919132954
$ws.Range("AK2").Value = 0.4263869195795579
$ws.Range("AN2").ClearContents()
$ws.Range("AP2").ClearContents()

# Row 3
$ws.Range("D3").ClearContents()
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = -37
$ws.Range("L3").Value = -1.002710027100271
$ws.Range("O3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("U3").Value = 285.4
$ws.Range("V3").Value = 2.162121212121212
$ws.Range("W3").Value = -0.06934032983508245
$ws.Range("X3").Value = 0.1336224303203382
$ws.Range("Y3").Value = -0.2029627601554206
$ws.Range("Z3").Value = 0.04997291440953412
$ws.Range("AA3").Value = 0
$ws.Range("AB3").Value = 0.06715997433322604
$ws.Range("AC3").Value = -0.06715997433322604
$ws.Range("AD3").Value = 375.3
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 375.3
$ws.Range("AG3").Value = 89.90000000000003
$ws.Range("AH3").Value = 0.7397989355411
$ws.Range("AI3").Value = 0.4131439894319683
$ws.Range("AJ3").Value = 0.4051374493014873
$ws.Range("AK3").Value = 0.1443017656500803
$ws.Range("AN3").ClearContents()
$ws.Range("AP3").ClearContents()

# Row 4
$ws.Range("B4").Value = "Eurobank Ergasias Services and Holdings S.A. (ATSE:EUROB)"
$ws.Range("D4").ClearContents()
$ws.Range("F4").Value = 0.257
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = -1230.6
$ws.Range("L4").Value = -8.545833333333333
$ws.Range("O4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("U4").Value = 7152.4
$ws.Range("V4").Value = 2.725450596349503
$ws.Range("W4").Value = -0.1722250990161365
$ws.Range("X4").Value = 0.2490925801383365
$ws.Range("Y4").Value = -0.4213176791544729
$ws.Range("Z4").Value = 0.007629745410231276
$ws.Range("AA4").Value = 0
$ws.Range("AB4").Value = 0.07042653807323118
$ws.Range("AC4").Value = -0.07042653807323118
$ws.Range("AD4").Value = 17598.6
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 17598.6
$ws.Range("AG4").Value = 10446.2
$ws.Range("AH4").Value = 0.8702312724683404
$ws.Range("AI4").Value = 0.7370060933475721
$ws.Range("AJ4").Value = 0.799219616694082
$ws.Range("AK4").Value = 0.6245448729829428
$ws.Range("AN4").ClearContents()
$ws.Range("AP4").ClearContents()

# Row 5
$ws.Range("B5").Value = "Alpha Bank A.E. (ATSE:ALPHA)"
$ws.Range("D5").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("F5").Value = 0.259
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 159.5
$ws.Range("L5").Value = 0.1055033734620982
$ws.Range("M5").Value = -0
$ws.Range("N5").Value = -0
$ws.Range("O5").Value = -0
$ws.Range("S5").Value = 0
$ws.Range("T5").ClearContents()
$ws.Range("U5").Value = 6126.7
$ws.Range("V5").Value = 3.400321900321901
$ws.Range("W5").Value = 0.0171544112111337
$ws.Range("X5").Value = 0.1332517177221639
$ws.Range("Y5").Value = -0.1160973065110302
$ws.Range("Z5").Value = 0.08164918610051958
$ws.Range("AA5").Value = 0
$ws.Range("AB5").Value = 0.08516646298621965
$ws.Range("AC5").Value = -0.08516646298621965
$ws.Range("AD5").Value = 5100.5
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 5100.5
$ws.Range("AG5").Value = -1026.2
$ws.Range("AH5").Value = 0.7389565796908277
$ws.Range("AI5").Value = 0.3392846452162229
$ws.Range("AJ5").Value = -1.323104693140794
$ws.Range("AK5").Value = -0.1152205155842989
$ws.Range("AN5").ClearContents()
$ws.Range("AP5").ClearContents()

# Row 6
$ws.Range("B6").Value = "National Bank of Greece S.A. (ATSE:ETE)"
$ws.Range("F6").ClearContents()
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = -203.1
$ws.Range("L6").Value = -0.1143903125880034
$ws.Range("M6").Value = 20
$ws.Range("N6").Value = 0.007911392405063292
$ws.Range("O6").Value = -0.09847365829640571
$ws.Range("R6").Value = 0
$ws.Range("S6").Value = 20
$ws.Range("U6").Value = 7926.3
$ws.Range("V6").Value = 3.135403481012658
$ws.Range("W6").Value = -0.03168289030325721
$ws.Range("X6").Value = 0.2680035396108901
$ws.Range("Y6").Value = -0.2996864299141473
$ws.Range("Z6").Value = 0.1343711686621157
$ws.Range("AA6").Value = 0
$ws.Range("AB6").Value = 0.09214344800337318
$ws.Range("AC6").Value = -0.09214344800337318
$ws.Range("AD6").Value = 18552.1
$ws.Range("AE6").Value = 0
$ws.Range("AF6").Value = 18552.1
$ws.Range("AG6").Value = 10625.8
$ws.Range("AH6").Value = 0.8800764702254733
$ws.Range("AI6").Value = 0.743657579438091
$ws.Range("AJ6").Value = 0.807812191153887
$ws.Range("AK6").Value = 0.6242832299304381
$ws.Range("AN6").ClearContents()
$ws.Range("AP6").ClearContents()

# Row 7
$ws.Range("F7").ClearContents()
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 49.3
$ws.Range("L7").Value = 0.03466947960618846
$ws.Range("M7").Value = 1.17
$ws.Range("N7").Value = 0.001686121919584954
$ws.Range("O7").Value = 0.02373225152129817
$ws.Range("S7").Value = 1.17
$ws.Range("T7").Value = 1
$ws.Range("U7").Value = 9336.5
$ws.Range("V7").Value = 13.45510880530336
$ws.Range("W7").Value = 0.005916378648233487
$ws.Range("X7").Value = 0.6120079809780188
$ws.Range("Y7").Value = -0.6060916023297853
$ws.Range("Z7").Value = 0.1317825865344516
$ws.Range("AA7").Value = 0
$ws.Range("AB7").Value = 0.1033758064503938
$ws.Range("AC7").Value = -0.1033758064503938
$ws.Range("AD7").Value = 13077.7
$ws.Range("AE7").Value = 0
$ws.Range("AF7").Value = 13077.7
$ws.Range("AG7").Value = 3741.200000000001
$ws.Range("AH7").Value = 0.9496136977547998
$ws.Range("AI7").Value = 0.5928671163819679
$ws.Range("AJ7").Value = 0.8435435503145364
$ws.Range("AK7").Value = 0.294075570472964
$ws.Range("AN7").ClearContents()
$ws.Range("AP7").ClearContents()
